# Applies scheduled-runner price/profit updates to the Brynhildr_Profits workbook.
# Each job-sheet (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR) gets refreshed market-board
# pricing (currentAveragePrice*, LevePrice*) and recomputed LeveProfit* columns.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 2860.1482
$ws.Range("I113").Value = 2632.2222
$ws.Range("J113").Value = 3316
$ws.Range("K113").Value = 2632.2222
$ws.Range("L113").Value = 3316
$ws.Range("M113").Value = 621.7777999999998
$ws.Range("N113").Value = -9824
$ws.Range("H130").Value = 275233
$ws.Range("J130").Value = 275233
$ws.Range("L130").Value = 275233
$ws.Range("N130").Value = -285273
$ws.Range("H132").Value = 10449.571
$ws.Range("I132").Value = 10449.571
$ws.Range("K132").Value = 31348.713
$ws.Range("M132").Value = -28818.713
$ws.Range("H138").Value = 4641.3
$ws.Range("I138").Value = 4641.3
$ws.Range("K138").Value = 13923.9
$ws.Range("M138").Value = -8783.900000000001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1163.5883
$ws.Range("I2").Value = 1163.8
$ws.Range("J2").Value = 1162
$ws.Range("K2").Value = 1163.8
$ws.Range("L2").Value = 1162
$ws.Range("M2").Value = -1050.8
$ws.Range("N2").Value = -1388
$ws.Range("H32").Value = 1071801.9
$ws.Range("J32").Value = 52000
$ws.Range("L32").Value = 52000
$ws.Range("N32").Value = -52574
$ws.Range("H45").Value = 2631.7058
$ws.Range("I45").Value = 1825.6666
$ws.Range("K45").Value = 1825.6666
$ws.Range("M45").Value = -1448.6666
$ws.Range("H116").Value = 1163.5883
$ws.Range("I116").Value = 1163.8
$ws.Range("J116").Value = 1162
$ws.Range("K116").Value = 1163.8
$ws.Range("L116").Value = 1162
$ws.Range("M116").Value = 1130.2
$ws.Range("N116").Value = -5750
$ws.Range("H122").Value = 3220.625
$ws.Range("I122").Value = 2791
$ws.Range("K122").Value = 8373
$ws.Range("M122").Value = -5923
$ws.Range("H128").Value = 30000
$ws.Range("J128").Value = 30000
$ws.Range("L128").Value = 30000
$ws.Range("N128").Value = -39960
$ws.Range("H132").Value = 5513.7407
$ws.Range("I132").Value = 4252.4546
$ws.Range("K132").Value = 12757.3638
$ws.Range("M132").Value = -10227.3638

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1163.5883
$ws.Range("I3").Value = 1163.8
$ws.Range("J3").Value = 1162
$ws.Range("K3").Value = 1163.8
$ws.Range("L3").Value = 1162
$ws.Range("M3").Value = -1049.8
$ws.Range("N3").Value = -1390
$ws.Range("H20").Value = 80813.71000000001
$ws.Range("J20").Value = 11080
$ws.Range("L20").Value = 11080
$ws.Range("N20").Value = -11574
$ws.Range("H76").Value = 18156.5
$ws.Range("J76").Value = 18156.5
$ws.Range("L76").Value = 18156.5
$ws.Range("N76").Value = -18786.5
$ws.Range("H79").Value = 18156.5
$ws.Range("J79").Value = 18156.5
$ws.Range("L79").Value = 18156.5
$ws.Range("N79").Value = -20340.5
$ws.Range("H80").Value = 470.27777
$ws.Range("I80").Value = 571.625
$ws.Range("K80").Value = 571.625
$ws.Range("M80").Value = 426.375
$ws.Range("H83").Value = 470.27777
$ws.Range("I83").Value = 571.625
$ws.Range("K83").Value = 2858.125
$ws.Range("M83").Value = 2133.875
$ws.Range("H86").Value = 2213.0625
$ws.Range("I86").Value = 2160.6
$ws.Range("K86").Value = 2160.6
$ws.Range("M86").Value = -1037.6
$ws.Range("H89").Value = 2213.0625
$ws.Range("I89").Value = 2160.6
$ws.Range("K89").Value = 10803
$ws.Range("M89").Value = -5187
$ws.Range("H99").Value = 6463.55
$ws.Range("I99").Value = 8380.214
$ws.Range("J99").Value = 1991.3334
$ws.Range("K99").Value = 8380.214
$ws.Range("L99").Value = 1991.3334
$ws.Range("M99").Value = -6882.214
$ws.Range("N99").Value = -4987.3334
$ws.Range("H134").Value = 9260741
$ws.Range("I134").Value = 1666.875
$ws.Range("K134").Value = 5000.625
$ws.Range("M134").Value = -2465.625

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10535015
$ws.Range("I31").Value = 10535015
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 10535015
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -10534720
$ws.Range("H34").Value = 10535015
$ws.Range("I34").Value = 10535015
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 10535015
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -10534813
$ws.Range("H82").Value = 5555
$ws.Range("I82").Value = 5555
$ws.Range("K82").Value = 5555
$ws.Range("M82").Value = -5194
$ws.Range("H85").Value = 5555
$ws.Range("I85").Value = 5555
$ws.Range("K85").Value = 5555
$ws.Range("M85").Value = -4307
$ws.Range("H99").Value = 25969.666
$ws.Range("I99").Value = 43709.6
$ws.Range("J99").Value = 3794.75
$ws.Range("K99").Value = 43709.6
$ws.Range("L99").Value = 3794.75
$ws.Range("M99").Value = -42211.6
$ws.Range("N99").Value = -6790.75
$ws.Range("H122").Value = 30373.666
$ws.Range("J122").Value = 65053.5
$ws.Range("L122").Value = 195160.5
$ws.Range("N122").Value = -200060.5
$ws.Range("H126").Value = 25969.666
$ws.Range("I126").Value = 43709.6
$ws.Range("J126").Value = 3794.75
$ws.Range("K126").Value = 131128.8
$ws.Range("L126").Value = 11384.25
$ws.Range("M126").Value = -128658.8
$ws.Range("N126").Value = -16324.25
$ws.Range("H132").Value = 5882.2964
$ws.Range("I132").Value = 6440.609
$ws.Range("J132").Value = 2672
$ws.Range("K132").Value = 19321.827
$ws.Range("L132").Value = 8016
$ws.Range("M132").Value = -16791.827
$ws.Range("N132").Value = -13076
$ws.Range("H134").Value = 2670.5366
$ws.Range("I134").Value = 2658.0303
$ws.Range("J134").Value = 2722.125
$ws.Range("K134").Value = 7974.090899999999
$ws.Range("L134").Value = 8166.375
$ws.Range("M134").Value = -5439.090899999999
$ws.Range("N134").Value = -13236.375
$ws.Range("H135").Value = 91284.86
$ws.Range("I135").Value = 40000
$ws.Range("J135").Value = 99832.336
$ws.Range("K135").Value = 40000
$ws.Range("L135").Value = 99832.336
$ws.Range("M135").Value = -34930
$ws.Range("N135").Value = -109972.336
$ws.Range("N31").ClearContents()
$ws.Range("N34").ClearContents()

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 46000
$ws.Range("J37").Value = 46000
$ws.Range("L37").Value = 138000
$ws.Range("N37").Value = -138224
$ws.Range("H107").Value = 3808.6
$ws.Range("I107").Value = 1297.4
$ws.Range("J107").Value = 5064.2
$ws.Range("K107").Value = 3892.2
$ws.Range("L107").Value = 15192.6
$ws.Range("M107").Value = -1972.2
$ws.Range("N107").Value = -19032.6
$ws.Range("H109").Value = 7220.9414
$ws.Range("I109").Value = 5551.143
$ws.Range("J109").Value = 8389.799999999999
$ws.Range("K109").Value = 16653.429
$ws.Range("L109").Value = 25169.4
$ws.Range("M109").Value = -15613.429
$ws.Range("N109").Value = -27249.4
$ws.Range("H113").Value = 573.7273
$ws.Range("J113").Value = 573.7273
$ws.Range("L113").Value = 1721.1819
$ws.Range("N113").Value = -6061.1819

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1254.375
$ws.Range("I102").Value = 1262.2858
$ws.Range("K102").Value = 1262.2858
$ws.Range("M102").Value = 359.7141999999999
$ws.Range("H126").Value = 3186
$ws.Range("I126").Value = 3107.5
$ws.Range("K126").Value = 9322.5
$ws.Range("M126").Value = -6852.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H69").Value = 56666.668
$ws.Range("I69").Value = 20000
$ws.Range("J69").Value = 75000
$ws.Range("K69").Value = 20000
$ws.Range("L69").Value = 75000
$ws.Range("M69").Value = -19189
$ws.Range("N69").Value = -76622
$ws.Range("H72").Value = 56666.668
$ws.Range("I72").Value = 20000
$ws.Range("J72").Value = 75000
$ws.Range("K72").Value = 60000
$ws.Range("L72").Value = 225000
$ws.Range("M72").Value = -55944
$ws.Range("N72").Value = -233112
$ws.Range("H99").Value = 44947
$ws.Range("I99").Value = 44947
$ws.Range("K99").Value = 44947
$ws.Range("M99").Value = -41952

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 12150.667
$ws.Range("I62").Value = 7594.625
$ws.Range("J62").Value = 17357.572
$ws.Range("K62").Value = 7594.625
$ws.Range("L62").Value = 17357.572
$ws.Range("M62").Value = -6970.625
$ws.Range("N62").Value = -18605.572
$ws.Range("H65").Value = 12150.667
$ws.Range("I65").Value = 7594.625
$ws.Range("J65").Value = 17357.572
$ws.Range("K65").Value = 37973.125
$ws.Range("L65").Value = 86787.86
$ws.Range("M65").Value = -34853.125
$ws.Range("N65").Value = -93027.86
$ws.Range("H81").Value = 80710.36
$ws.Range("J81").Value = 277424.75
$ws.Range("L81").Value = 554849.5
$ws.Range("N81").Value = -556971.5
$ws.Range("H82").Value = 30364.166
$ws.Range("J82").Value = 31237
$ws.Range("L82").Value = 31237
$ws.Range("N82").Value = -32003
$ws.Range("H84").Value = 80710.36
$ws.Range("J84").Value = 277424.75
$ws.Range("L84").Value = 2774247.5
$ws.Range("N84").Value = -2784855.5
$ws.Range("H85").Value = 30364.166
$ws.Range("J85").Value = 31237
$ws.Range("L85").Value = 31237
$ws.Range("N85").Value = -33889
$ws.Range("H96").Value = 2924
$ws.Range("I96").Value = 2849
$ws.Range("K96").Value = 2849
$ws.Range("M96").Value = -1476
$ws.Range("H126").Value = 1496
$ws.Range("I126").Value = 1186.75
$ws.Range("J126").Value = 2733
$ws.Range("K126").Value = 3560.25
$ws.Range("L126").Value = 8199
$ws.Range("M126").Value = -1090.25
$ws.Range("N126").Value = -13139
$ws.Range("H136").Value = 8229853
$ws.Range("I136").Value = 4348689
$ws.Range("J136").Value = 21167066
$ws.Range("K136").Value = 13046067
$ws.Range("L136").Value = 63501198
$ws.Range("M136").Value = -13043517
$ws.Range("N136").Value = -63506298
